$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$r = $p.Range
$markOnly = $d.Range($r.End - 1, $r.End)
$markOnly.Delete()
